$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# F1 header: copy style from E1 (bold/centered/bordered) then set text
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)
$ws.Range("F1").Value = "time_taken"

# F2:F115 - time_taken metadata values (plain, unstyled data cells)
$timeTaken = @(
    "2021-10-05 10:51:57.981055",
    "2021-10-05 10:51:57.981067",
    "2021-10-05 10:51:57.981070",
    "2021-10-05 10:51:57.981073",
    "2021-10-05 10:51:57.981075",
    "2021-10-05 10:51:57.981078",
    "2021-10-05 10:51:57.981081",
    "2021-10-05 10:51:57.981083",
    "2021-10-05 10:51:57.981086",
    "2021-10-05 10:51:57.981088",
    "2021-10-05 10:51:57.981091",
    "2021-10-05 10:51:57.981093",
    "2021-10-05 10:51:57.981096",
    "2021-10-05 10:51:57.981098",
    "2021-10-05 10:51:57.981101",
    "2021-10-05 10:51:57.981103",
    "2021-10-05 10:51:57.981106",
    "2021-10-05 10:51:57.981108",
    "2021-10-05 10:51:57.981111",
    "2021-10-05 10:51:57.981113",
    "2021-10-05 10:51:57.981116",
    "2021-10-05 10:51:57.981118",
    "2021-10-05 10:51:57.981121",
    "2021-10-05 10:51:57.981124",
    "2021-10-05 10:51:57.981127",
    "2021-10-05 10:51:57.981129",
    "2021-10-05 10:51:57.981132",
    "2021-10-05 10:51:57.981134",
    "2021-10-05 10:51:57.981137",
    "2021-10-05 10:51:57.981139",
    "2021-10-05 10:51:57.981142",
    "2021-10-05 10:51:57.981144",
    "2021-10-05 10:51:57.981147",
    "2021-10-05 10:51:57.981150",
    "2021-10-05 10:51:57.981152",
    "2021-10-05 10:51:57.981155",
    "2021-10-05 10:51:57.981157",
    "2021-10-05 10:51:57.981160",
    "2021-10-05 10:51:57.981162",
    "2021-10-05 10:51:57.981165",
    "2021-10-05 10:51:57.981168",
    "2021-10-05 10:51:57.981170",
    "2021-10-05 10:51:57.981173",
    "2021-10-05 10:51:57.981175",
    "2021-10-05 10:51:57.981178",
    "2021-10-05 10:51:57.981180",
    "2021-10-05 10:51:57.981183",
    "2021-10-05 10:51:57.981185",
    "2021-10-05 10:51:57.981188",
    "2021-10-05 10:51:57.981190",
    "2021-10-05 10:51:57.981193",
    "2021-10-05 10:51:57.981195",
    "2021-10-05 10:51:57.981198",
    "2021-10-05 10:51:57.981201",
    "2021-10-05 10:51:57.981203",
    "2021-10-05 10:51:57.981206",
    "2021-10-05 10:51:57.981208",
    "2021-10-05 10:51:57.981211",
    "2021-10-05 10:51:57.981213",
    "2021-10-05 10:51:57.981216",
    "2021-10-05 10:51:57.981218",
    "2021-10-05 10:51:57.981221",
    "2021-10-05 10:51:57.981223",
    "2021-10-05 10:51:57.981226",
    "2021-10-05 10:51:57.981229",
    "2021-10-05 10:51:57.981232",
    "2021-10-05 10:51:57.981235",
    "2021-10-05 10:51:57.981237",
    "2021-10-05 10:51:57.981240",
    "2021-10-05 10:51:57.981242",
    "2021-10-05 10:51:57.981245",
    "2021-10-05 10:51:57.981247",
    "2021-10-05 10:51:57.981249",
    "2021-10-05 10:51:57.981252",
    "2021-10-05 10:51:57.981254",
    "2021-10-05 10:51:57.981257",
    "2021-10-05 10:51:57.981261",
    "2021-10-05 10:51:57.981264",
    "2021-10-05 10:51:57.981267",
    "2021-10-05 10:51:57.981269",
    "2021-10-05 10:51:57.981272",
    "2021-10-05 10:51:57.981274",
    "2021-10-05 10:51:57.981277",
    "2021-10-05 10:51:57.981279",
    "2021-10-05 10:51:57.981282",
    "2021-10-05 10:51:57.981284",
    "2021-10-05 10:51:57.981287",
    "2021-10-05 10:51:57.981289",
    "2021-10-05 10:51:57.981292",
    "2021-10-05 10:51:57.981294",
    "2021-10-05 10:51:57.981297",
    "2021-10-05 10:51:57.981300",
    "2021-10-05 10:51:57.981303",
    "2021-10-05 10:51:57.981306",
    "2021-10-05 10:51:57.981308",
    "2021-10-05 10:51:57.981311",
    "2021-10-05 10:51:57.981313",
    "2021-10-05 10:51:57.981316",
    "2021-10-05 10:51:57.981318",
    "2021-10-05 10:51:57.981321",
    "2021-10-05 10:51:57.981323",
    "2021-10-05 10:51:57.981326",
    "2021-10-05 10:51:57.981328",
    "2021-10-05 10:51:57.981331",
    "2021-10-05 10:51:57.981333",
    "2021-10-05 10:51:57.981336",
    "2021-10-05 10:51:57.981338",
    "2021-10-05 10:51:57.981341",
    "2021-10-05 10:51:57.981345",
    "2021-10-05 10:51:57.981348",
    "2021-10-05 10:51:57.981351",
    "2021-10-05 10:51:57.981353",
    "2021-10-05 10:51:57.981356",
    "2021-10-05 10:51:57.981358"
)

for ($i = 0; $i -lt $timeTaken.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 6).Value = $timeTaken[$i]
}
